$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = '64.258.76'
$c.ClearFormats()
$c = $ws.Cells.Item(2, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.44%  '
$c.ClearFormats()
$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = '3.146.93'
$c.ClearFormats()
$c = $ws.Cells.Item(3, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.83%  '
$c.ClearFormats()
$c = $ws.Cells.Item(4, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.02%  '
$c.ClearFormats()
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = '591.27'
$c.ClearFormats()
$c = $ws.Cells.Item(5, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.82%  '
$c.ClearFormats()
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = '145.89'
$c.ClearFormats()
$c = $ws.Cells.Item(6, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.29%  '
$c.ClearFormats()
$c = $ws.Cells.Item(7, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.01%  '
$c.ClearFormats()
$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = '3.138.55'
$c.ClearFormats()
$c = $ws.Cells.Item(8, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.73%  '
$c.ClearFormats()
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = '0.530'
$c.ClearFormats()
$c = $ws.Cells.Item(9, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.15%  '
$c.ClearFormats()
$c = $ws.Cells.Item(10, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.28%  '
$c.ClearFormats()
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = '5.96'
$c.ClearFormats()
$c = $ws.Cells.Item(11, 5)
$c.NumberFormat = "@"
$c.Value = '  +3.09%  '
$c.ClearFormats()
$c = $ws.Cells.Item(12, 5)
$c.NumberFormat = "@"
$c.Value = '  -1.12%  '
$c.ClearFormats()
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = '0.0000248'
$c.ClearFormats()
$c = $ws.Cells.Item(13, 5)
$c.NumberFormat = "@"
$c.Value = '  -1.41%  '
$c.ClearFormats()
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = '37.44'
$c.ClearFormats()
$c = $ws.Cells.Item(14, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.47%  '
$c.ClearFormats()
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = '3.667.14'
$c.ClearFormats()
$c = $ws.Cells.Item(15, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.87%  '
$c.ClearFormats()
$c = $ws.Cells.Item(16, 5)
$c.NumberFormat = "@"
$c.Value = '  -1.17%  '
$c.ClearFormats()
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = '7.32'
$c.ClearFormats()
$c = $ws.Cells.Item(17, 5)
$c.NumberFormat = "@"
$c.Value = '  +2.41%  '
$c.ClearFormats()
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = '64.027.23'
$c.ClearFormats()
$c = $ws.Cells.Item(18, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.27%  '
$c.ClearFormats()
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = '3.145.49'
$c.ClearFormats()
$c = $ws.Cells.Item(19, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.73%  '
$c.ClearFormats()
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = '470.07'
$c.ClearFormats()
$c = $ws.Cells.Item(20, 5)
$c.NumberFormat = "@"
$c.Value = '  +1.11%  '
$c.ClearFormats()
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = '14.38'
$c.ClearFormats()
$c = $ws.Cells.Item(21, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.27%  '
$c.ClearFormats()
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = '0.734'
$c.ClearFormats()
$c = $ws.Cells.Item(22, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.47%  '
$c.ClearFormats()
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = '7.59'
$c.ClearFormats()
$c = $ws.Cells.Item(23, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.71%  '
$c.ClearFormats()
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = '2.37'
$c.ClearFormats()
$c = $ws.Cells.Item(24, 5)
$c.NumberFormat = "@"
$c.Value = '  +9.39%  '
$c.ClearFormats()
$c = $ws.Cells.Item(25, 5)
$c.NumberFormat = "@"
$c.Value = '  -1.42%  '
$c.ClearFormats()
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = '81.49'
$c.ClearFormats()
$c = $ws.Cells.Item(26, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.36%  '
$c.ClearFormats()
$c = $ws.Cells.Item(27, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.16%  '
$c.ClearFormats()
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = '9.99'
$c.ClearFormats()
$c = $ws.Cells.Item(28, 5)
$c.NumberFormat = "@"
$c.Value = '  +11.91%  '
$c.ClearFormats()
$c = $ws.Cells.Item(30, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.89%  '
$c.ClearFormats()
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = '2.24'
$c.ClearFormats()
$c = $ws.Cells.Item(31, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.48%  '
$c.ClearFormats()
$c = $ws.Cells.Item(32, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.15%  '
$c.ClearFormats()
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = '27.68'
$c.ClearFormats()
$c = $ws.Cells.Item(33, 5)
$c.NumberFormat = "@"
$c.Value = '  +2.46%  '
$c.ClearFormats()
$c = $ws.Cells.Item(34, 5)
$c.NumberFormat = "@"
$c.Value = '  +1.15%  '
$c.ClearFormats()
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = '0.0₃0847'
$c.ClearFormats()
$c = $ws.Cells.Item(35, 5)
$c.NumberFormat = "@"
$c.Value = '  -2.38%  '
$c.ClearFormats()
$c = $ws.Cells.Item(36, 5)
$c.NumberFormat = "@"
$c.Value = '  +1.16%  '
$c.ClearFormats()
$c = $ws.Cells.Item(37, 2)
$c.NumberFormat = "@"
$c.Value = 'Stacks'
$c.ClearFormats()
$c = $ws.Cells.Item(37, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c.ClearFormats()
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = '2.32'
$c.ClearFormats()
$c = $ws.Cells.Item(37, 5)
$c.NumberFormat = "@"
$c.Value = '  -1.96%  '
$c.ClearFormats()
$c = $ws.Cells.Item(38, 2)
$c.NumberFormat = "@"
$c.Value = 'Filecoin'
$c.ClearFormats()
$c = $ws.Cells.Item(38, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c.ClearFormats()
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = '6.16'
$c.ClearFormats()
$c = $ws.Cells.Item(38, 5)
$c.NumberFormat = "@"
$c.Value = '  +1.66%  '
$c.ClearFormats()
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = '3.22'
$c.ClearFormats()
$c = $ws.Cells.Item(39, 5)
$c.NumberFormat = "@"
$c.Value = '  -5.22%  '
$c.ClearFormats()
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = '51.48'
$c.ClearFormats()
$c = $ws.Cells.Item(40, 5)
$c.NumberFormat = "@"
$c.Value = '  +1.02%  '
$c.ClearFormats()
$c = $ws.Cells.Item(41, 5)
$c.NumberFormat = "@"
$c.Value = '  +6.95%  '
$c.ClearFormats()
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = '455.81'
$c.ClearFormats()
$c = $ws.Cells.Item(42, 5)
$c.NumberFormat = "@"
$c.Value = '  +1.66%  '
$c.ClearFormats()
$c = $ws.Cells.Item(43, 5)
$c.NumberFormat = "@"
$c.Value = '  +6.55%  '
$c.ClearFormats()
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = '0.0374'
$c.ClearFormats()
$c = $ws.Cells.Item(44, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.57%  '
$c.ClearFormats()
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = '2.925.89'
$c.ClearFormats()
$c = $ws.Cells.Item(45, 5)
$c.NumberFormat = "@"
$c.Value = '  +1.62%  '
$c.ClearFormats()
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = '40.66'
$c.ClearFormats()
$c = $ws.Cells.Item(46, 5)
$c.NumberFormat = "@"
$c.Value = '  +13.63%  '
$c.ClearFormats()
$c = $ws.Cells.Item(47, 5)
$c.NumberFormat = "@"
$c.Value = '  -2.50%  '
$c.ClearFormats()
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = '133.99'
$c.ClearFormats()
$c = $ws.Cells.Item(48, 5)
$c.NumberFormat = "@"
$c.Value = '  +8.46%  '
$c.ClearFormats()
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = '2.25'
$c.ClearFormats()
$c = $ws.Cells.Item(50, 5)
$c.NumberFormat = "@"
$c.Value = '  +2.97%  '
$c.ClearFormats()
$c = $ws.Cells.Item(51, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.31%  '
$c.ClearFormats()
